# Generate Report for handback
# Refresh the "Latest Handoff Datetime" and "Latest Handback DateTime" for the
# first (16be212c...) file row on each localized-language report sheet.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("D2").Value = "2016-01-08 07:20:45"
$ws_zhcn.Range("G2").Value = "2016-01-08 07:21:31"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("D2").Value = "2016-01-08 07:20:55"
$ws_dede.Range("G2").Value = "2016-01-08 07:21:50"
